$wb = $excel.ActiveWorkbook

# --- Sheet1: add header row A1:H1 ---------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$values1 = @("id", "name", "item", "kaka", "kaka", "item", "item", "item")
for ($i = 0; $i -lt $values1.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $values1[$i]
}

# --- Sheet33 (3rd sheet): add header row A1:B1 ---------------------------
$ws3 = $wb.Worksheets.Item(3)
$values3 = @("id", "name")
for ($i = 0; $i -lt $values3.Length; $i++) {
    $ws3.Cells.Item(1, $i + 1).Value = $values3[$i]
}

# --- Make Sheet1 the active sheet/tab, with G4 selected -------------------
[void]$ws1.Activate()
[void]$ws1.Range("G4").Select()

Write-Output "ok"
